$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dados Base")

# Remove the "DepreciacaoAmortizacao" row (row 11), shifting all subsequent
# rows up by one and shrinking the used range from A1:F17 to A1:F16.
$ws.Rows.Item(11).Delete()
